$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values to match the final content/layout ---
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "LOB1010"
$ws.Range("C2").Value = "LOB1010"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Direito  Aplicado à Engenharia"
$ws.Range("C3").Value = " Direito  Aplicado à Engenharia"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Law applied to engineering"
$ws.Range("C4").Value = "Law applied to engineering"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2012"
$ws.Range("C8").Value = "01/01/2012"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EF-9,EM-8,EB-8,EP-10,EQD-8,EQN-11"
$ws.Range("C9").Value = "EF-9,EM-8,EB-8,EP-10,EQD-8,EQN-11"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C10").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C18").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Provas"
$ws.Range("C19").Value = "Provas"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "NF= (P1+P2)/2"
$ws.Range("C20").Value = "NF= (P1+P2)/2"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Reestudo com trabalhos e prova"
$ws.Range("C21").Value = "Reestudo com trabalhos e prova"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n"
$ws.Range("C23").Value = "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n"

# --- Clear cells that existed before but have no content after ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Row heights: set explicit heights, then delete the now-removed row 24 ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30

# Rows that must go back to default (no explicit height) -> AutoFit removes custom height flag
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()

# --- Remove the now-deleted last row (24) to shrink used range / dimension to A1:C23 ---
$ws.Rows.Item(24).Delete()
